$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = -4163

function Set-TextValue($ws, $addr, $val) {
    # Write via a text formula then paste-special as values so Excel
    # stores the literal text (e.g. "1.00", "48.314.69") without
    # reinterpreting number-like strings as numeric values.
    $escaped = $val.Replace('"', '""')
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

Set-TextValue $ws "D2" "48.314.69"
Set-TextValue $ws "E2" "  +2.07%  "
Set-TextValue $ws "D3" "2.506.12"
Set-TextValue $ws "E3" "  +0.65%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.03%  "
Set-TextValue $ws "D5" "321.58"
Set-TextValue $ws "E5" "  -0.03%  "
Set-TextValue $ws "D6" "108.34"
Set-TextValue $ws "E6" "  -0.74%  "
Set-TextValue $ws "D7" "0.528"
Set-TextValue $ws "E7" "  +1.15%  "
Set-TextValue $ws "D8" "1.00"
Set-TextValue $ws "E8" "  +0.03%  "
Set-TextValue $ws "D9" "0.542"
Set-TextValue $ws "D10" "39.99"
Set-TextValue $ws "E10" "  +1.44%  "
Set-TextValue $ws "E11" "  +9.46%  "
Set-TextValue $ws "D12" "0.0819"
Set-TextValue $ws "E12" "  +1.19%  "
Set-TextValue $ws "E13" "  +0.11%  "
Set-TextValue $ws "D14" "7.20"
Set-TextValue $ws "E14" "  +0.26%  "
Set-TextValue $ws "D15" "2.898.32"
Set-TextValue $ws "E15" "  +0.67%  "
Set-TextValue $ws "D16" "2.502.16"
Set-TextValue $ws "E16" "  +0.43%  "
Set-TextValue $ws "D17" "0.847"
Set-TextValue $ws "E17" "  +0.24%  "
Set-TextValue $ws "D18" "48.156.14"
Set-TextValue $ws "E18" "  +1.92%  "
Set-TextValue $ws "D19" "13.12"
Set-TextValue $ws "E19" "  -2.19%  "
Set-TextValue $ws "D20" "6.77"
Set-TextValue $ws "E20" "  +2.17%  "
Set-TextValue $ws "E21" "  +0.80%  "
Set-TextValue $ws "E22" "  +1.60%  "
Set-TextValue $ws "D23" "280.80"
Set-TextValue $ws "E23" "  +13.81%  "
Set-TextValue $ws "D24" "72.14"
Set-TextValue $ws "E24" "  +2.19%  "
Set-TextValue $ws "D25" "2.56"
Set-TextValue $ws "E25" "  +0.35%  "
Set-TextValue $ws "E26" "  +0.02%  "
Set-TextValue $ws "E27" "  +0.43%  "
Set-TextValue $ws "D28" "2.26"
Set-TextValue $ws "E28" "  -1.43%  "
Set-TextValue $ws "D29" "9.81"
Set-TextValue $ws "E29" "  -1.54%  "
Set-TextValue $ws "D30" "0.141"
Set-TextValue $ws "E30" "  +1.18%  "
Set-TextValue $ws "D31" "35.34"
Set-TextValue $ws "E31" "  +2.17%  "
Set-TextValue $ws "D32" "49.33"
Set-TextValue $ws "E32" "  -1.11%  "
Set-TextValue $ws "D33" "19.58"
Set-TextValue $ws "E33" "  -4.22%  "
Set-TextValue $ws "D34" "5.37"
Set-TextValue $ws "E34" "  +1.17%  "
Set-TextValue $ws "E35" "  -0.08%  "
Set-TextValue $ws "D36" "0.0786"
Set-TextValue $ws "E36" "  +0.19%  "
Set-TextValue $ws "E37" "  +0.00%  "
Set-TextValue $ws "D38" "4.67"
Set-TextValue $ws "E38" "  -1.43%  "
Set-TextValue $ws "E39" "  -0.04%  "
Set-TextValue $ws "E40" "  +0.12%  "
Set-TextValue $ws "D41" "121.50"
Set-TextValue $ws "E41" "  +1.75%  "
Set-TextValue $ws "E42" "  +0.41%  "
Set-TextValue $ws "D43" "21.58"
Set-TextValue $ws "E43" "  -4.77%  "
Set-TextValue $ws "D44" "0.0305"
Set-TextValue $ws "E44" "  +2.84%  "
Set-TextValue $ws "D45" "2.012.47"
Set-TextValue $ws "E45" "  +0.90%  "
Set-TextValue $ws "D46" "3.20"
Set-TextValue $ws "E46" "  +5.62%  "
Set-TextValue $ws "D47" "1.86"
Set-TextValue $ws "E47" "  +4.36%  "
Set-TextValue $ws "E48" "  -2.13%  "
Set-TextValue $ws "D49" "8.98"
Set-TextValue $ws "E49" "  -1.44%  "
Set-TextValue $ws "D50" "5.19"
Set-TextValue $ws "E50" "  +0.26%  "
Set-TextValue $ws "D51" "80.39"
Set-TextValue $ws "E51" "  +3.77%  "
